# Edit script: applies the diff changes to the document via Word COM interop.
$d = $word.ActiveDocument

# --- 1. "Around 1" + "2" + "+" + " years in the retail domain" -> "Around 12+ years in the retail domain"
$d.Content.Find.Execute("Around 12+ years in the retail domain", $false, $false, $false, $false, $false, $true, 1, $false, "Around 12+ years in the retail domain", 2)

Write-Host "step1 done"

# --- 2. "Supplier " + "Collaboration platform (B2B)" -> "Supplier Collaboration platform (B2B)"
$d.Content.Find.Execute("Supplier Collaboration platform (B2B)", $false, $false, $false, $false, $false, $true, 1, $false, "Supplier Collaboration platform (B2B)", 2)

Write-Host "step2 done"

# --- 3. "•" + " " + "Contributed to system analysis, process mapping, and " -> "• Contributed to system analysis, process mapping, and "
$d.Content.Find.Execute([char]0x2022 + " Contributed to system analysis, process mapping, and ", $false, $false, $false, $false, $false, $true, 1, $false, [char]0x2022 + " Contributed to system analysis, process mapping, and ", 2)

Write-Host "step3 done"

# --- 4. "•" + " " + "Delivered a modernization roadmap to streamline workflows and " -> "• Delivered a modernization roadmap to streamline workflows and "
$d.Content.Find.Execute([char]0x2022 + " Delivered a modernization roadmap to streamline workflows and ", $false, $false, $false, $false, $false, $true, 1, $false, [char]0x2022 + " Delivered a modernization roadmap to streamline workflows and ", 2)

Write-Host "step4 done"

# --- 5. Awards table: resize columns/indent, re-center header cols 2 & 3,
#        un-center data-row col 1, and shorten "Star of the Month award" text.
$awardsTable = $d.Tables.Item(2)

# Overall table width + indent (dxa -> points, 20 dxa per point)
$awardsTable.PreferredWidth = 508.5

# Row indent (tblInd) for every row
for ($r = 1; $r -le $awardsTable.Rows.Count; $r++) {
    $awardsTable.Rows.Item($r).LeftIndent = -27.25
}

# Column widths (dxa -> points)
$awardsTable.Columns.Item(1).Width = 99.0
$awardsTable.Columns.Item(2).Width = 252.0
$awardsTable.Columns.Item(3).Width = 94.5

# Header row (row 1): center-align columns 2 and 3
$headerRow = $awardsTable.Rows.Item(1)
$headerRow.Cells.Item(2).Range.Paragraphs.Item(1).Alignment = 1
$headerRow.Cells.Item(3).Range.Paragraphs.Item(1).Alignment = 1

# Data rows (rows 2-9): remove center-alignment from column 1
for ($r = 2; $r -le $awardsTable.Rows.Count; $r++) {
    $awardsTable.Rows.Item($r).Cells.Item(1).Range.Paragraphs.Item(1).Alignment = 0
}

Write-Host "step5a done"

# --- 6. "Star of the Month award" -> "Star of the Month" (row 2, col 1)
$d.Content.Find.Execute("Star of the Month award", $true, $false, $false, $false, $false, $true, 1, $false, "Star of the Month", 2)

Write-Host "step5b done"
